$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.141.09'
$ws.Range('E2').Value = '  +2.86%  '
$ws.Range('D3').Value = '2.301.14'
$ws.Range('E3').Value = '  +1.76%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.28'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.538'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.41%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.525'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +7.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.07'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0821'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.98%  '
$ws.Range('E12').Value = '  +0.99%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.15'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.38%  '
$ws.Range('D14').Value = '2.652.53'
$ws.Range('E14').Value = '  +1.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.05'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.55%  '
$ws.Range('D16').Value = '2.298.76'
$ws.Range('E16').Value = '  +1.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.810'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.42%  '
$ws.Range('D18').Value = '43.030.88'
$ws.Range('E18').Value = '  +2.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.76%  '
$ws.Range('D20').Value = '0.0₃0926'
$ws.Range('E20').Value = '  +2.81%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.08'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.49'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.60'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.53%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.63'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.82%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.01'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.94%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.61'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.44'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.32'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +10.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.64'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.97%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '167.49'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.31'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.13'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.78'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0741'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.82%  '
$ws.Range('E37').Value = '  +2.49%  '
$ws.Range('E39').Value = '  +1.79%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.83'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('E41').Value = '  +5.65%  '
$ws.Range('E42').Value = '  -0.51%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0289'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.35%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.974.85'
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.01'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.05'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.85'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('E48').Value = '  +17.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '55.58'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.53%  '
$ws.Range('D50').Value = '2.522.79'
$ws.Range('E50').Value = '  +1.55%  '
$ws.Range('E51').Value = '  +2.54%  '
